$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A1").ClearContents()
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Value = "是是是"
}
[void]$ws.Range("A2:A13").Select()
